# Weekly update: add this week's "Fruta" (Frambuesa) record at the top of
# the log (row 8, just under the header block) and push the previously
# existing weekly rows (old rows 8-18) down by one row (new rows 9-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (old rows 8-18) down by inserting a new
# blank row at row 8; Excel's native row-insert semantics take care of
# moving all the old row 8..18 content down to rows 9..19.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with this week's entry. All of the "dimension"
# style columns (Mercado/Región/Producto/etc.) repeat unchanged from the
# rest of the table; only the date, quality grade, volume and the four
# price columns are new for this week.
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = 'Vega Monumental Concepción'
$ws.Range("C8").Value = 'Bíobío'
$ws.Range("D8").Value = 44944
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = 'Berries'
$ws.Range("I8").Value = 100101004
$ws.Range("J8").Value = 'Frambuesa'
$ws.Range("K8").Value = 'Sin especificar'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 7500
$ws.Range("Q8").Value = '$/bandeja 2 kilos'
$ws.Range("R8").Value = 'Región de Ñuble'
$ws.Range("S8").Value = 3750
$ws.Range("T8").Value = 2
